# daily auto push: 2026-01-17 13:37 UTC
#
# A new observation (2026/01/17, Saturday "土", hour 19) needs to be
# inserted into the time-series table on Sheet1. It belongs right before
# the existing "2026/12/29" block, at row 657, so every row from the old
# 657 down to the old 698 shifts down by one (dimension grows from
# A1:D698 to A1:D699).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 657..698 down one slot, leaving a blank row 657 to populate.
$ws.Rows("657:657").Insert()

# Columns A (date) and B (weekday) hold text that looks like a date /
# single kanji character. Mark them as Text before assigning so Excel's
# autodetect doesn't silently reinterpret "2026/01/17" as a date serial.
# ClearFormats() afterwards drops that temporary number format again so
# the new cells end up with the same (default) style as their neighbours.
$ws.Range("A657:B657").NumberFormat = "@"
$ws.Range("A657").Value = "2026/01/17"
$ws.Range("B657").Value = "土"
$ws.Range("C657").Value = 19
$ws.Range("D657").Value = 201
$ws.Range("A657:B657").ClearFormats()
